$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A12").Value = "IF (ROLE_SUBORDINATE IS high) THEN (SPEED_VOICE IS high_speed)"
$ws.Range("B12").Value = "chat_qualifier"

$ws.Range("A13").Value = "IF (ROLE_SUBORDINATE IS low) THEN (SPEED_VOICE IS mid_speed)"
$ws.Range("B13").Value = "chat_qualifier"

$ws.Range("A14").Value = "IF (ROLE_SUBORDINATE IS high) THEN (VOLUME IS high_volume)"
$ws.Range("B14").Value = "chat_qualifier"

$ws.Range("A15").Select()
